$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-09-19 -> 2023-09-20, serial 45188 -> 45189) for every data row (2..426).
$lastRow = 426
$range = $ws.Range("C2:C" + $lastRow)
$range.Value = 45189
